# DOMA-11339: add "Decommissioning date" column (Q) to the property meter
# import example sheet, with a sample value "2022-01-25" on the first data row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy the formatting (style/number format/borders) of column P - the current
# last column - onto the new column Q so the new cells look consistent with
# the rest of the header/data table (same text format, borders, etc).
$ws.Range("P1:P11").Copy()
$ws.Range("Q1:Q11").PasteSpecial(-4122)

# Give column Q the same width as column P.
$ws.Columns("Q").ColumnWidth = $ws.Columns("P").ColumnWidth

# Header for the new column.
$ws.Range("Q1").Value = "Decommissioning date"

# Sample decommissioning date for the first data row only (rest stay empty,
# same as the other example rows in the template).
$ws.Range("Q2").Value = "2022-01-25"

"Added Decommissioning date column (Q) with example value on row 2"
